$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data rows (row 2 .. row 6), columns A..J
# A = trial index, B..I = numeric schedule values, J = "train_dim2_1" (shared string)
$data = @(
    @(1, 3, 7, 5, 3, 2, -4, 32, 5),
    @(2, 2, 7, 3, 2, 1, -5, 21, 5),
    @(3, 1, 8, 6, 7, 5, -1, 65, 5),
    @(4, 4, 9, 8, 7, 4, -2, 54, 5),
    @(5, 2, 5, 5, 2, 3, -3, 43, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

# Update the selection shown in the sheet view to I1
$ws.Range("I1").Select()
